# Skill Course Weighting Master
# - Remove the stale AutoFilter criteria on the edges sheet (un-hides all
#   previously filtered-out rows) and drop the now-irrelevant "BA 450
#   Bachelor-Abschlussarbeit / Scientific Work (SW)" edge row, which shifts
#   every following row up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Clear the filter criteria and turn the AutoFilter off entirely - this
# un-hides every row that the old filter had hidden.
$ws.ShowAllData()
$ws.AutoFilterMode = $false

# Delete the "BA 450 Bachelor-Abschlussarbeit Bachelor Thesis" edge row;
# everything below shifts up by one row.
$ws.Rows(83).Delete()

# Re-apply a clean AutoFilter (no active filter criteria) over the new
# data extent.
$ws.UsedRange.AutoFilter()

# Restore the user's on-screen selection to the (new) last row.
$ws.Activate()
$ws.Rows(83).Select()
